$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update product name (B1) on both sheets - shared text value
$wsInput.Range("B1").Value = "4134-RBI-EI-DB-SAR-NOREC-MOREREPAY-1st"
$wsOutput.Range("B1").Value = "4134-RBI-EI-DB-SAR-NOREC-MOREREPAY-1st"

# Update short name (B2) on input sheet - now a text value instead of numeric
$wsInput.Range("B2").Value = "413y"

# Update the selected cell on the input sheet to B3, then restore the originally active tab
$wsInput.Range("B3").Select()
$wsOutput.Activate()
